$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '45.696.53'
Set-TextValue $ws.Range('E2') '  +6.12%  '
Set-TextValue $ws.Range('D3') '2.400.34'
Set-TextValue $ws.Range('E3') '  +4.19%  '
Set-TextValue $ws.Range('E4') '  +0.08%  '
Set-TextValue $ws.Range('D5') '115.47'
Set-TextValue $ws.Range('E5') '  +9.97%  '
Set-TextValue $ws.Range('D6') '320.04'
Set-TextValue $ws.Range('E6') '  +2.92%  '
Set-TextValue $ws.Range('E7') '  +2.15%  '
Set-TextValue $ws.Range('E8') '  -0.16%  '
Set-TextValue $ws.Range('D9') '0.631'
Set-TextValue $ws.Range('E9') '  +4.04%  '
Set-TextValue $ws.Range('D10') '42.89'
Set-TextValue $ws.Range('E10') '  +7.85%  '
Set-TextValue $ws.Range('E11') '  +3.16%  '
Set-TextValue $ws.Range('D12') '8.73'
Set-TextValue $ws.Range('E12') '  +5.25%  '
Set-TextValue $ws.Range('E13') '  +3.11%  '
Set-TextValue $ws.Range('E14') '  +2.95%  '
Set-TextValue $ws.Range('D15') '16.04'
Set-TextValue $ws.Range('E15') '  +4.26%  '
Set-TextValue $ws.Range('D16') '2.766.69'
Set-TextValue $ws.Range('E16') '  -0.59%  '
Set-TextValue $ws.Range('D17') '2.403.61'
Set-TextValue $ws.Range('E17') '  +4.46%  '
Set-TextValue $ws.Range('D18') '45.713.74'
Set-TextValue $ws.Range('E18') '  +6.57%  '
Set-TextValue $ws.Range('D19') '7.54'
Set-TextValue $ws.Range('E19') '  +2.68%  '
Set-TextValue $ws.Range('E20') '  +3.56%  '
Set-TextValue $ws.Range('D21') '13.55'
Set-TextValue $ws.Range('E21') '  -0.37%  '
Set-TextValue $ws.Range('D22') '75.28'
Set-TextValue $ws.Range('E22') '  +2.50%  '
Set-TextValue $ws.Range('D23') '3.59'
Set-TextValue $ws.Range('E23') '  +4.06%  '
Set-TextValue $ws.Range('D24') '266.29'
Set-TextValue $ws.Range('E24') '  -0.64%  '
Set-TextValue $ws.Range('D25') '2.38'
Set-TextValue $ws.Range('E25') '  +7.63%  '
Set-TextValue $ws.Range('D26') '0.999'
Set-TextValue $ws.Range('E26') '  -0.75%  '
Set-TextValue $ws.Range('D27') '7.67'
Set-TextValue $ws.Range('E27') '  +5.37%  '
Set-TextValue $ws.Range('E28') '  +4.58%  '
Set-TextValue $ws.Range('E29') '  +2.94%  '
Set-TextValue $ws.Range('D30') '40.44'
Set-TextValue $ws.Range('E30') '  +11.35%  '
Set-TextValue $ws.Range('E31') '  +16.53%  '
Set-TextValue $ws.Range('E32') '  +2.50%  '
Set-TextValue $ws.Range('D33') '173.19'
Set-TextValue $ws.Range('E33') '  +5.01%  '
Set-TextValue $ws.Range('E34') '  +12.24%  '
Set-TextValue $ws.Range('B35') 'RenderToken'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D35') '5.05'
Set-TextValue $ws.Range('E35') '  +11.08%  '
Set-TextValue $ws.Range('B36') 'Stellar'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D36') '0.133'
Set-TextValue $ws.Range('E36') '  +1.89%  '
Set-TextValue $ws.Range('D37') '0.120'
Set-TextValue $ws.Range('E37') '  +7.76%  '
Set-TextValue $ws.Range('D38') '4.24'
Set-TextValue $ws.Range('E38') '  +16.89%  '
Set-TextValue $ws.Range('D39') '3.12'
Set-TextValue $ws.Range('E39') '  +10.19%  '
Set-TextValue $ws.Range('D40') '0.0366'
Set-TextValue $ws.Range('E40') '  +5.28%  '
Set-TextValue $ws.Range('E41') '  +11.22%  '
Set-TextValue $ws.Range('D42') '0.245'
Set-TextValue $ws.Range('E42') '  +7.80%  '
Set-TextValue $ws.Range('D43') '13.84'
Set-TextValue $ws.Range('E43') '  +12.49%  '
Set-TextValue $ws.Range('D44') '100.68'
Set-TextValue $ws.Range('E44') '  -8.83%  '
Set-TextValue $ws.Range('D45') '72.51'
Set-TextValue $ws.Range('E45') '  +2.14%  '
Set-TextValue $ws.Range('D46') '90.91'
Set-TextValue $ws.Range('E46') '  +16.70%  '
Set-TextValue $ws.Range('E47') '  -0.44%  '
Set-TextValue $ws.Range('D48') '5.89'
Set-TextValue $ws.Range('E48') '  +14.23%  '
Set-TextValue $ws.Range('D49') '116.80'
Set-TextValue $ws.Range('E49') '  +5.29%  '
Set-TextValue $ws.Range('E50') '  +9.80%  '
Set-TextValue $ws.Range('D51') '1.59'
Set-TextValue $ws.Range('E51') '  +11.58%  '
